$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 983.5714
$ws.Range("J32").Value = 1017.4
$ws.Range("L32").Value = 1017.4
$ws.Range("N32").Value = -1669.4
$ws.Range("H38").Value = 14435.571
$ws.Range("I38").Value = 190
$ws.Range("J38").Value = 50049.5
$ws.Range("K38").Value = 570
$ws.Range("L38").Value = 150148.5
$ws.Range("M38").Value = -198
$ws.Range("N38").Value = -150892.5
$ws.Range("H42").Value = 53.666668
$ws.Range("J42").Value = 21
$ws.Range("L42").Value = 63
$ws.Range("N42").Value = -523
$ws.Range("H86").Value = 27166.334
$ws.Range("J86").Value = 1500
$ws.Range("L86").Value = 1500
$ws.Range("N86").Value = -3746
$ws.Range("H89").Value = 27166.334
$ws.Range("J89").Value = 1500
$ws.Range("L89").Value = 7500
$ws.Range("N89").Value = -18732
$ws.Range("H100").Value = 3261.3333
$ws.Range("I100").Value = 1740
$ws.Range("K100").Value = 1740
$ws.Range("M100").Value = -1199
$ws.Range("H112").Value = 2370.76
$ws.Range("J112").Value = 2274.8096
$ws.Range("L112").Value = 6824.4288
$ws.Range("N112").Value = -9040.4288
$ws.Range("H113").Value = 4079.1177
$ws.Range("I113").Value = 3054.6667
$ws.Range("J113").Value = 5231.625
$ws.Range("K113").Value = 3054.6667
$ws.Range("L113").Value = 5231.625
$ws.Range("M113").Value = 199.3332999999998
$ws.Range("N113").Value = -11739.625
$ws.Range("H125").Value = 3982
$ws.Range("J125").Value = 5284.8
$ws.Range("L125").Value = 47563.2
$ws.Range("N125").Value = -52483.2
$ws.Range("H138").Value = 2695.01
$ws.Range("J138").Value = 2925.5625
$ws.Range("L138").Value = 8776.6875
$ws.Range("N138").Value = -19056.6875

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2011
$ws.Range("I2").Value = 1164.1428
$ws.Range("K2").Value = 1164.1428
$ws.Range("M2").Value = -1051.1428
$ws.Range("H45").Value = 4571.567
$ws.Range("I45").Value = 4981.636
$ws.Range("J45").Value = 4334.1577
$ws.Range("K45").Value = 4981.636
$ws.Range("L45").Value = 4334.1577
$ws.Range("M45").Value = -4604.636
$ws.Range("N45").Value = -5088.1577
$ws.Range("H63").Value = 2993.7778
$ws.Range("I63").Value = 1573.5
$ws.Range("K63").Value = 1573.5
$ws.Range("M63").Value = -887.5
$ws.Range("H66").Value = 2993.7778
$ws.Range("I66").Value = 1573.5
$ws.Range("K66").Value = 7867.5
$ws.Range("M66").Value = -4435.5
$ws.Range("H102").Value = 1412.1578
$ws.Range("I102").Value = 1166.5883
$ws.Range("K102").Value = 1166.5883
$ws.Range("M102").Value = 455.4117000000001
$ws.Range("H116").Value = 2011
$ws.Range("I116").Value = 1164.1428
$ws.Range("K116").Value = 1164.1428
$ws.Range("M116").Value = 1129.8572
$ws.Range("H132").Value = 4445.4585
$ws.Range("I132").Value = 4927.5
$ws.Range("J132").Value = 3770.6
$ws.Range("K132").Value = 14782.5
$ws.Range("L132").Value = 11311.8
$ws.Range("M132").Value = -12252.5
$ws.Range("N132").Value = -16371.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2011
$ws.Range("I3").Value = 1164.1428
$ws.Range("K3").Value = 1164.1428
$ws.Range("M3").Value = -1050.1428
$ws.Range("H11").Value = 749
$ws.Range("I11").Value = 198.33333
$ws.Range("J11").Value = 1575
$ws.Range("K11").Value = 198.33333
$ws.Range("L11").Value = 1575
$ws.Range("M11").Value = -58.33332999999999
$ws.Range("N11").Value = -1855
$ws.Range("H12").Value = 1698.3334
$ws.Range("I12").Value = 100
$ws.Range("K12").Value = 100
$ws.Range("M12").Value = 68
$ws.Range("H117").Value = 119000
$ws.Range("J117").Value = 119000
$ws.Range("L117").Value = 119000
$ws.Range("N117").Value = -128178

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2343.9092
$ws.Range("J22").Value = 1196
$ws.Range("L22").Value = 1196
$ws.Range("N22").Value = -1896
$ws.Range("H31").Value = 6461
$ws.Range("J31").Value = 7916.875
$ws.Range("L31").Value = 7916.875
$ws.Range("N31").Value = -8506.875
$ws.Range("H34").Value = 6461
$ws.Range("J34").Value = 7916.875
$ws.Range("L34").Value = 7916.875
$ws.Range("N34").Value = -8320.875
$ws.Range("H99").Value = 2453
$ws.Range("I99").Value = 2453
$ws.Range("K99").Value = 2453
$ws.Range("M99").Value = -955
$ws.Range("H122").Value = 5274.231
$ws.Range("I122").Value = 4368.5713
$ws.Range("J122").Value = 6330.8335
$ws.Range("K122").Value = 13105.7139
$ws.Range("L122").Value = 18992.5005
$ws.Range("M122").Value = -10655.7139
$ws.Range("N122").Value = -23892.5005
$ws.Range("H126").Value = 2453
$ws.Range("I126").Value = 2453
$ws.Range("K126").Value = 7359
$ws.Range("M126").Value = -4889

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 929.6667
$ws.Range("I5").Value = 894.5
$ws.Range("K5").Value = 2683.5
$ws.Range("M5").Value = -2571.5
$ws.Range("H59").Value = 1599.6
$ws.Range("J59").Value = 2033
$ws.Range("L59").Value = 6099
$ws.Range("N59").Value = -7179
$ws.Range("H81").Value = 3518.5
$ws.Range("I81").Value = 1111
$ws.Range("J81").Value = 4000
$ws.Range("K81").Value = 3333
$ws.Range("L81").Value = 12000
$ws.Range("M81").Value = -2210
$ws.Range("N81").Value = -14246
$ws.Range("H84").Value = 3518.5
$ws.Range("I84").Value = 1111
$ws.Range("J84").Value = 4000
$ws.Range("K84").Value = 9999
$ws.Range("L84").Value = 36000
$ws.Range("M84").Value = -4383
$ws.Range("N84").Value = -47232
$ws.Range("H109").Value = 1156.5
$ws.Range("I109").Value = 1156.5
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 3469.5
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -2429.5
$ws.Range("N109").Value = $null
$ws.Range("H113").Value = 2145.0833
$ws.Range("I113").Value = 1096
$ws.Range("K113").Value = 3288
$ws.Range("M113").Value = -1118
$ws.Range("H131").Value = 1561.5116
$ws.Range("J131").Value = 1742.9445
$ws.Range("L131").Value = 5228.833500000001
$ws.Range("N131").Value = -15308.8335
$ws.Range("H135").Value = 929.6667
$ws.Range("I135").Value = 894.5
$ws.Range("K135").Value = 8050.5
$ws.Range("M135").Value = -5515.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 810.2857
$ws.Range("I107").Value = 810.2857
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 810.2857
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1109.7143
$ws.Range("N107").Value = $null
$ws.Range("H113").Value = 25518.586
$ws.Range("I113").Value = 10002.5
$ws.Range("J113").Value = 99995.8
$ws.Range("K113").Value = 10002.5
$ws.Range("L113").Value = 99995.8
$ws.Range("M113").Value = -7832.5
$ws.Range("N113").Value = -104335.8
$ws.Range("H122").Value = 2635.75
$ws.Range("I122").Value = 2635.75
$ws.Range("K122").Value = 7907.25
$ws.Range("M122").Value = -5457.25
$ws.Range("H126").Value = 2722
$ws.Range("I126").Value = 2517.2
$ws.Range("K126").Value = 7551.599999999999
$ws.Range("M126").Value = -5081.599999999999
$ws.Range("H132").Value = 3999.75
$ws.Range("I132").Value = 4250
$ws.Range("J132").Value = 3749.5
$ws.Range("K132").Value = 12750
$ws.Range("L132").Value = 11248.5
$ws.Range("M132").Value = -10220
$ws.Range("N132").Value = -16308.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 15126.934
$ws.Range("J122").Value = 11981
$ws.Range("L122").Value = 35943
$ws.Range("N122").Value = -40843

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4719.3228
$ws.Range("I122").Value = 5076.2607
$ws.Range("J122").Value = 3693.125
$ws.Range("K122").Value = 15228.7821
$ws.Range("L122").Value = 11079.375
$ws.Range("M122").Value = -12778.7821
$ws.Range("N122").Value = -15979.375
$ws.Range("H126").Value = 3549.75
$ws.Range("I126").Value = 2959.6
$ws.Range("J126").Value = 4533.3335
$ws.Range("K126").Value = 8878.799999999999
$ws.Range("L126").Value = 13600.0005
$ws.Range("M126").Value = -6408.799999999999
$ws.Range("N126").Value = -18540.0005
$ws.Range("H136").Value = 2655.5264
$ws.Range("I136").Value = 2135.8667
$ws.Range("K136").Value = 6407.6001
$ws.Range("M136").Value = -3857.6001
